$wb = $excel.ActiveWorkbook

# Sheet ALC Row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 975.62
$ws.Range("I15").Value = 975.62
$ws.Range("K15").Value = 2926.86
$ws.Range("M15").Value = -2757.86

# Sheet ALC Row 87
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 15192.835
$ws.Range("J87").Value = 15192.835
$ws.Range("L87").Value = 15192.835
$ws.Range("N87").Value = -17688.835

# Sheet ALC Row 90
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 15192.835
$ws.Range("J90").Value = 15192.835
$ws.Range("L90").Value = 45578.505
$ws.Range("N90").Value = -58058.505

# Sheet ALC Row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1218.6086
$ws.Range("J112").Value = 1271.4
$ws.Range("L112").Value = 3814.2
$ws.Range("N112").Value = -6030.200000000001

# Sheet ALC Row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 2778.5
$ws.Range("I116").Value = 2540.7144
$ws.Range("J116").Value = 3333.3333
$ws.Range("K116").Value = 2540.7144
$ws.Range("L116").Value = 3333.3333
$ws.Range("M116").Value = 901.2856000000002
$ws.Range("N116").Value = -10217.3333

# Sheet ALC Row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3697.2444
$ws.Range("I138").Value = 2074.258
$ws.Range("J138").Value = 4550
$ws.Range("K138").Value = 6222.773999999999
$ws.Range("L138").Value = 13650
$ws.Range("M138").Value = -1082.773999999999
$ws.Range("N138").Value = -23930

# Sheet ARM Row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2679.8333
$ws.Range("I2").Value = 2679.8333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2679.8333
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -2566.8333

# Sheet ARM Row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1288.3334
$ws.Range("I45").Value = 1072.5
$ws.Range("J45").Value = 1720
$ws.Range("K45").Value = 1072.5
$ws.Range("L45").Value = 1720
$ws.Range("M45").Value = -695.5
$ws.Range("N45").Value = -2474

# Sheet ARM Row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2272.4849
$ws.Range("I61").Value = 2122.182
$ws.Range("J61").Value = 2573.0908
$ws.Range("K61").Value = 2122.182
$ws.Range("L61").Value = 2573.0908
$ws.Range("M61").Value = -1910.182
$ws.Range("N61").Value = -2997.0908

# Sheet ARM Row 62
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 23699.75
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 29933
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 29933
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -31181

# Sheet ARM Row 65
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H65").Value = 23699.75
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 29933
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 89799
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -96039

# Sheet ARM Row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 2679.8333
$ws.Range("I116").Value = 2679.8333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 2679.8333
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -385.8332999999998

# Sheet ARM Row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 2272.4849
$ws.Range("I136").Value = 2122.182
$ws.Range("J136").Value = 2573.0908
$ws.Range("K136").Value = 6366.545999999999
$ws.Range("L136").Value = 7719.2724
$ws.Range("M136").Value = -3816.545999999999
$ws.Range("N136").Value = -12819.2724

# Sheet BSM Row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2679.8333
$ws.Range("I3").Value = 2679.8333
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2679.8333
$ws.Range("L3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -2565.8333

# Sheet BSM Row 50
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H50").Value = 36345
$ws.Range("J50").Value = 36345
$ws.Range("L50").Value = 36345
$ws.Range("N50").Value = -37493

# Sheet CRP Row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3832.5403
$ws.Range("I31").Value = 1548.9246
$ws.Range("K31").Value = 1548.9246
$ws.Range("M31").Value = -1253.9246

# Sheet CRP Row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3832.5403
$ws.Range("I34").Value = 1548.9246
$ws.Range("K34").Value = 1548.9246
$ws.Range("M34").Value = -1346.9246

# Sheet CRP Row 64
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 30000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 30000
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").Value = 30000
$ws.Range("N64").Value = -30496

# Sheet CRP Row 67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H67").Value = 30000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 30000
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").Value = 30000
$ws.Range("N67").Value = -31716

# Sheet CRP Row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1318817.6
$ws.Range("I132").Value = 2501998.8
$ws.Range("J132").Value = 4172
$ws.Range("K132").Value = 7505996.399999999
$ws.Range("L132").Value = 12516
$ws.Range("M132").Value = -7503466.399999999
$ws.Range("N132").Value = -17576

# Sheet CRP Row 141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 38536.125
$ws.Range("I141").Value = 19334.334
$ws.Range("J141").Value = 42967.31
$ws.Range("K141").Value = 19334.334
$ws.Range("L141").Value = 42967.31
$ws.Range("M141").Value = -14154.334
$ws.Range("N141").Value = -53327.31

# Sheet CUL Row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 1064.7142
$ws.Range("I98").Value = 1130.6
$ws.Range("J98").Value = 900
$ws.Range("K98").Value = 3391.8
$ws.Range("L98").Value = 2700
$ws.Range("M98").Value = -1893.8
$ws.Range("N98").Value = -5696

# Sheet CUL Row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 4325.2856
$ws.Range("I99").Value = 1649.5
$ws.Range("J99").Value = 5395.6
$ws.Range("K99").Value = 4948.5
$ws.Range("L99").Value = 16186.8
$ws.Range("M99").Value = -2702.5
$ws.Range("N99").Value = -20678.8

# Sheet CUL Row 100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 5169.478
$ws.Range("J100").Value = 5313.5454
$ws.Range("L100").Value = 15940.6362
$ws.Range("N100").Value = -17562.6362

# Sheet CUL Row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 842.5
$ws.Range("I131").Value = 200
$ws.Range("J131").Value = 869.2708
$ws.Range("K131").Value = 600
$ws.Range("L131").Value = 2607.8124
$ws.Range("M131").Value = 4440
$ws.Range("N131").Value = -12687.8124

# Sheet CUL Row 132
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1611.8889
$ws.Range("I132").Value = 1104
$ws.Range("J132").Value = 1757
$ws.Range("K132").Value = 9936
$ws.Range("L132").Value = 15813
$ws.Range("M132").Value = -7406
$ws.Range("N132").Value = -20873

# Sheet LTW Row 68
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 3277.1428
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3490
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 3490
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -4988

# Sheet LTW Row 71
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 3277.1428
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3490
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 17450
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -24938

# Sheet WVR Row 40
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 17000
$ws.Range("I40").Value = 17000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 17000
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -16851

# Sheet WVR Row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2840
$ws.Range("I62").Value = 2980
$ws.Range("J62").Value = 2490
$ws.Range("K62").Value = 2980
$ws.Range("L62").Value = 2490
$ws.Range("M62").Value = -2356
$ws.Range("N62").Value = -3738

# Sheet WVR Row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 2840
$ws.Range("I65").Value = 2980
$ws.Range("J65").Value = 2490
$ws.Range("K65").Value = 14900
$ws.Range("L65").Value = 12450
$ws.Range("M65").Value = -11780
$ws.Range("N65").Value = -18690

# Sheet WVR Row 138
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H138").Value = 50976.332
$ws.Range("J138").Value = 50976.332
$ws.Range("L138").Value = 50976.332
$ws.Range("N138").Value = -61256.332
